$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the data block that had 11 rows instead of 10: clear row 276's content
# (old A276=0,B276=0), keeping subsequent row numbers intact (no shifting).
$ws.Range("A276:B276").ClearContents()

# Complete the last block (label row 288) to 10 data rows by adding row 298
$ws.Range("A298").Value = 0
$ws.Range("B298").Value = 0

# Add synthese (average) columns D/E for each block of 10 data rows.
$firstRows = @(2,13,24,35,46,57,68,79,90,101,112,123,134,145,156,167,178,189,200,211,222,233,244,255,266,278,289)

foreach ($r in $firstRows) {
    $endR = $r + 9
    $ws.Range("D" + $r).Formula = "=AVERAGE(A" + $r + ":A" + $endR + ")"
    $ws.Range("E" + $r).Formula = "=AVERAGE(B" + $r + ":B" + $endR + ")"
}

$ws.Range("A280").Select()
$ws.Range("D278").Select()
